$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status active -> draft
$ws.Range("B6").Value = "draft"

# Date updated
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact row 1 (row 10): new publisher-style contact text
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact row 2 (row 11): now has the named contact instead of the old duplicate text
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row 12 for "Jurisdiction" (pushes old rows 12-15 down to 13-16),
# copying formatting from the row below so borders/fill match the rest of the table.
$ws.Rows.Item(12).Insert()
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
